$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing cell values
$ws.Range("L3").Value = 5
$ws.Range("F5").Value = -8

# Append a new row (row 10) duplicating row 8's data
$ws.Range("A10").Value = 733688164476661
$ws.Range("B10").Value = 5630279
$ws.Range("C10").Value = "F"
$ws.Range("D10").Value = "2016-04-27T15:05:12Z"
$ws.Range("E10").Value = "2016-04-29T00:00:00Z"
$ws.Range("F10").Value = 23
$ws.Range("G10").Value = "GOIABEIRAS"
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = 0
$ws.Range("N10").Value = "Yes"

# Update the active selection to match the saved view state
$ws.Range("P9").Select()
